$d = $word.ActiveDocument

# Locate "finger….which" (ellipsis + period between "finger" and "which")
# so we can compute character offsets without hard-coding absolute
# positions.
$ellipsis = [char]0x2026
$findRange = $d.Content
$findRange.Find.Execute("finger" + $ellipsis + ".which", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null

$matchStart = $findRange.Start          # start of "finger"
$afterFinge = $matchStart + 5           # between "finge" and "r"
$afterR     = $matchStart + 6           # between "r" and the ellipsis/period
$beforeWhich = $matchStart + 8          # after the ellipsis+period, before "which"

# 1) Replace the ellipsis + period with "; "
$rngPunct = $d.Range($afterR, $beforeWhich)
$rngPunct.Text = "; "

# 2) Force the run to split between "finge" and "r; " using a throwaway
#    bookmark (adding/removing a bookmark splits runs without leaving any
#    residual formatting behind).
$splitRng = $d.Range($afterFinge, $afterFinge)
$d.Bookmarks.Add("TempSplit", $splitRng)
$d.Bookmarks("TempSplit").Delete()

# 3) Relocate the "_GoBack" bookmark so it sits right after "r; " and
#    before "which", matching where Word leaves it after this edit.
$goBackRng = $d.Range($beforeWhich, $beforeWhich)
$d.Bookmarks.Add("_GoBack", $goBackRng)
